$d = $word.ActiveDocument

# The "Chuyen giao san pham" (Product delivery) table is the 4th table
# in the document: Deliverable Name / Description / Delivery Date.
$t = $d.Tables.Item(4)

function Replace-InCell($table, $row, $col, $findText, $replaceText) {
    # Re-resolve the cell range fresh (offsets can shift after edits) and
    # rebuild it through $d.Range(...) so Find.Execute stays confined to
    # the cell instead of drifting into other parts of the document.
    $cellRange = $table.Cell($row, $col).Range
    $scoped = $d.Range($cellRange.Start, $cellRange.End)
    $scoped.Find.ClearFormatting()
    $result = $scoped.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 0, $false, $replaceText, 2)
    return $result
}

# Row 5, col 1: "Chức năng tập đọc,  luyện từ và câu, kể chuyện"
#            -> "Chức năng tập đọc,  luyện từ và câu."
Replace-InCell $t 5 1 ", kể chuyện" "." | Out-Null

# Row 5, col 2: "... tập đọc, luyện từ và câu, kể chuyện và các tài liệu ..."
#            -> "... tập đọc, luyện từ và câu và các tài liệu ..."
Replace-InCell $t 5 2 ", luyện từ và câu, kể chuyện" ", luyện từ và câu" | Out-Null

# Row 6, col 1: "Chức năng tập làm văn, viết nhật ký, game"
#            -> "Chức năng kể chuyện, tập làm văn."
Replace-InCell $t 6 1 "Chức năng tập làm văn, viết nhật ký, game" "Chức năng kể chuyện, tập làm văn." | Out-Null

# Row 7, col 1: "Chức năng đố vui, kể chuyện cổ tích"
#            -> "Chức năng viết nhật ký, game."
Replace-InCell $t 7 1 "Chức năng đố vui, kể chuyện cổ tích" "Chức năng viết nhật ký, game." | Out-Null

# Append a new deliverable row: Sản phẩm hoàn chỉnh / Đã test xong. / 13/06
$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "Sản phẩm hoàn chỉnh"
$newRow.Cells.Item(2).Range.Text = "Đã test xong."
$newRow.Cells.Item(3).Range.Text = "13/06"
